# Daily attendance processing - 2025-11-15 09:43:06
# Normalizes the "Recorded By" (column G) entries on the Session Analysis
# Results sheet: whenever a cell holds a comma-separated list of recorders
# that does NOT already start with "System", the list is reversed so that
# "System" (when present) is surfaced first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -eq $null) { continue }
    if ($value -eq "") { continue }

    $parts = $value -split ", "
    $count = $parts.Count

    if ($count -le 1) { continue }
    if ($parts[0] -eq "System") { continue }

    $reversed = @()
    for ($i = $count - 1; $i -ge 0; $i--) {
        $reversed += $parts[$i]
    }

    $cell.Value = $reversed -join ", "
}
